$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Cells.Item(2, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '26.994.20'
$ws.Cells.Item(2, 5).Value = "  -0.51%  "

$dCell = $ws.Cells.Item(3, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.827.54'
$ws.Cells.Item(3, 5).Value = "  +0.21%  "

$ws.Cells.Item(4, 5).Value = "  -0.57%  "

$dCell = $ws.Cells.Item(5, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '312.27'
$ws.Cells.Item(5, 5).Value = "  -0.04%  "

$dCell = $ws.Cells.Item(6, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.005'
$ws.Cells.Item(6, 5).Value = "  -0.40%  "

$dCell = $ws.Cells.Item(7, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.4590'
$ws.Cells.Item(7, 5).Value = "  -0.76%  "

$dCell = $ws.Cells.Item(8, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.3699'
$ws.Cells.Item(8, 5).Value = "  +2.07%  "

$dCell = $ws.Cells.Item(9, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.07321'
$ws.Cells.Item(9, 5).Value = "  +0.43%  "

$dCell = $ws.Cells.Item(10, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.8729'
$ws.Cells.Item(10, 5).Value = "  +0.48%  "

$dCell = $ws.Cells.Item(11, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.07948'
$ws.Cells.Item(11, 5).Value = "  +3.99%  "

$dCell = $ws.Cells.Item(12, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '19.77'
$ws.Cells.Item(12, 5).Value = "  -1.60%  "

$dCell = $ws.Cells.Item(13, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.829.24'
$ws.Cells.Item(13, 5).Value = "  -0.25%  "

$dCell = $ws.Cells.Item(14, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '6.560'
$ws.Cells.Item(14, 5).Value = "  +1.44%  "

$dCell = $ws.Cells.Item(15, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '5.325'
$ws.Cells.Item(15, 5).Value = "  -0.22%  "

$dCell = $ws.Cells.Item(16, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '91.50'
$ws.Cells.Item(16, 5).Value = "  -0.95%  "

$dCell = $ws.Cells.Item(17, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.008'
$ws.Cells.Item(17, 5).Value = "  -0.25%  "

$dCell = $ws.Cells.Item(18, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.000008920'
$ws.Cells.Item(18, 5).Value = "  +3.43%  "

$dCell = $ws.Cells.Item(19, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.006'
$ws.Cells.Item(19, 5).Value = "  -0.46%  "

$dCell = $ws.Cells.Item(20, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '14.68'
$ws.Cells.Item(20, 5).Value = "  +1.58%  "

$dCell = $ws.Cells.Item(21, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '27.167.87'
$ws.Cells.Item(21, 5).Value = "  -1.00%  "

$dCell = $ws.Cells.Item(22, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '5.093'
$ws.Cells.Item(22, 5).Value = "  -2.20%  "

$dCell = $ws.Cells.Item(23, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '10.54'
$ws.Cells.Item(23, 5).Value = "  -0.03%  "

$dCell = $ws.Cells.Item(24, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '2.003.49'
$ws.Cells.Item(24, 5).Value = "  -4.45%  "

$dCell = $ws.Cells.Item(25, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '153.31'
$ws.Cells.Item(25, 5).Value = "  +1.50%  "

$dCell = $ws.Cells.Item(26, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.845'
$ws.Cells.Item(26, 5).Value = "  -1.87%  "

$dCell = $ws.Cells.Item(27, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '18.36'
$ws.Cells.Item(27, 5).Value = "  +1.02%  "

$ws.Cells.Item(28, 5).Value = "  -1.42%  "

$dCell = $ws.Cells.Item(29, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '5.137'
$ws.Cells.Item(29, 5).Value = "  +0.91%  "

$dCell = $ws.Cells.Item(30, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '114.85'
$ws.Cells.Item(30, 5).Value = "  -1.04%  "

$dCell = $ws.Cells.Item(31, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.08862'
$ws.Cells.Item(31, 5).Value = "  -0.40%  "

$dCell = $ws.Cells.Item(32, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '2.957'
$ws.Cells.Item(32, 5).Value = "  -0.19%  "

$dCell = $ws.Cells.Item(33, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.7308'
$ws.Cells.Item(33, 5).Value = "  -0.57%  "

$dCell = $ws.Cells.Item(34, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '4.442'
$ws.Cells.Item(34, 5).Value = "  -0.14%  "

$dCell = $ws.Cells.Item(35, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.130'
$ws.Cells.Item(35, 5).Value = "  -0.11%  "

$dCell = $ws.Cells.Item(36, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.072'
$ws.Cells.Item(36, 5).Value = "  -0.72%  "

$ws.Cells.Item(37, 5).Value = "  +1.76%  "

$ws.Cells.Item(38, 2).Value = 'Hedera'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$dCell = $ws.Cells.Item(38, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.05233'
$ws.Cells.Item(38, 5).Value = "  -0.08%  "

$ws.Cells.Item(39, 2).Value = 'RenderToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$dCell = $ws.Cells.Item(39, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '2.429'
$ws.Cells.Item(39, 5).Value = "  -2.44%  "

$dCell = $ws.Cells.Item(40, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '2.942'
$ws.Cells.Item(40, 5).Value = "  +0.51%  "

$dCell = $ws.Cells.Item(41, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '7.151'
$ws.Cells.Item(41, 5).Value = "  -0.04%  "

$dCell = $ws.Cells.Item(42, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.5141'
$ws.Cells.Item(42, 5).Value = "  -0.88%  "

$dCell = $ws.Cells.Item(43, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.1629'
$ws.Cells.Item(43, 5).Value = "  +0.36%  "

$dCell = $ws.Cells.Item(44, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '8.207'
$ws.Cells.Item(44, 5).Value = "  -0.85%  "

$dCell = $ws.Cells.Item(45, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.4828'
$ws.Cells.Item(45, 5).Value = "  +0.11%  "

$ws.Cells.Item(46, 5).Value = "  -0.45%  "

$dCell = $ws.Cells.Item(47, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '10.22'
$ws.Cells.Item(47, 5).Value = "  +0.73%  "

$dCell = $ws.Cells.Item(48, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '102.40'
$ws.Cells.Item(48, 5).Value = "  -0.91%  "

$dCell = $ws.Cells.Item(49, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.624'
$ws.Cells.Item(49, 5).Value = "  -0.48%  "

$dCell = $ws.Cells.Item(50, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.06211'
$ws.Cells.Item(50, 5).Value = "  -0.91%  "

$dCell = $ws.Cells.Item(51, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '64.56'
$ws.Cells.Item(51, 5).Value = "  +0.22%  "
